# Auto-generated Excel COM-interop script
# Applies numeric corrections to the Sagittarius_Profits workbook sheets
# (profit/price recalculations) as captured in the source diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 3375.2354  # H98
$ws.Cells.Item(98, 9).Value = 737.3  # I98
$ws.Cells.Item(98, 11).Value = 737.3  # K98
$ws.Cells.Item(98, 13).Value = 760.7  # M98

$ws.Cells.Item(107, 8).Value = 510.7143  # H107
$ws.Cells.Item(107, 9).Value = 575.4  # I107
$ws.Cells.Item(107, 11).Value = 575.4  # K107
$ws.Cells.Item(107, 13).Value = 1344.6  # M107

$ws.Cells.Item(112, 8).Value = 1545.75  # H112
$ws.Cells.Item(112, 9).Value = 1249  # I112
$ws.Cells.Item(112, 11).Value = 3747  # K112
$ws.Cells.Item(112, 13).Value = -2639  # M112

$ws.Cells.Item(122, 8).Value = 3375.2354  # H122
$ws.Cells.Item(122, 9).Value = 737.3  # I122
$ws.Cells.Item(122, 11).Value = 2211.9  # K122
$ws.Cells.Item(122, 13).Value = 238.1000000000004  # M122

$ws.Cells.Item(132, 8).Value = 1057.7  # H132
$ws.Cells.Item(132, 9).Value = 1064.1666  # I132
$ws.Cells.Item(132, 11).Value = 3192.4998  # K132
$ws.Cells.Item(132, 13).Value = -662.4998000000001  # M132

$ws.Cells.Item(138, 8).Value = 3162.5056  # H138
$ws.Cells.Item(138, 10).Value = 3114  # J138
$ws.Cells.Item(138, 12).Value = 9342  # L138
$ws.Cells.Item(138, 14).Value = -19622  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5588.9062  # H32
$ws.Cells.Item(32, 9).Value = 5422.4287  # I32
$ws.Cells.Item(32, 10).Value = 6754.25  # J32
$ws.Cells.Item(32, 11).Value = 5422.4287  # K32
$ws.Cells.Item(32, 12).Value = 6754.25  # L32
$ws.Cells.Item(32, 13).Value = -5135.4287  # M32
$ws.Cells.Item(32, 14).Value = -7328.25  # N32

$ws.Cells.Item(61, 8).Value = 4992.75  # H61
$ws.Cells.Item(61, 9).Value = 1986.6666  # I61
$ws.Cells.Item(61, 11).Value = 1986.6666  # K61
$ws.Cells.Item(61, 13).Value = -1774.6666  # M61

$ws.Cells.Item(97, 8).Value = 876.3333  # H97
$ws.Cells.Item(97, 9).Value = 911.6  # I97
$ws.Cells.Item(97, 11).Value = 911.6  # K97
$ws.Cells.Item(97, 13).Value = -415.6  # M97

$ws.Cells.Item(110, 8).Value = 1295.2858  # H110
$ws.Cells.Item(110, 9).Value = 1121.2727  # I110
$ws.Cells.Item(110, 11).Value = 1121.2727  # K110
$ws.Cells.Item(110, 13).Value = 923.7273  # M110

$ws.Cells.Item(122, 8).Value = 935.6  # H122
$ws.Cells.Item(122, 9).Value = 837.4737  # I122
$ws.Cells.Item(122, 11).Value = 2512.4211  # K122
$ws.Cells.Item(122, 13).Value = -62.42110000000002  # M122

$ws.Cells.Item(132, 8).Value = 1920  # H132
$ws.Cells.Item(132, 9).Value = 1682  # I132
$ws.Cells.Item(132, 11).Value = 5046  # K132
$ws.Cells.Item(132, 13).Value = -2516  # M132

$ws.Cells.Item(136, 8).Value = 4992.75  # H136
$ws.Cells.Item(136, 9).Value = 1986.6666  # I136
$ws.Cells.Item(136, 11).Value = 5959.9998  # K136
$ws.Cells.Item(136, 13).Value = -3409.9998  # M136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2797.3333  # H86
$ws.Cells.Item(86, 9).Value = 2502  # I86
$ws.Cells.Item(86, 11).Value = 2502  # K86
$ws.Cells.Item(86, 13).Value = -1379  # M86

$ws.Cells.Item(89, 8).Value = 2797.3333  # H89
$ws.Cells.Item(89, 9).Value = 2502  # I89
$ws.Cells.Item(89, 11).Value = 12510  # K89
$ws.Cells.Item(89, 13).Value = -6894  # M89

$ws.Cells.Item(94, 8).Value = 11668.167  # H94
$ws.Cells.Item(94, 9).Value = 11668.167  # I94
$ws.Cells.Item(94, 11).Value = 11668.167  # K94
$ws.Cells.Item(94, 13).Value = -11217.167  # M94

$ws.Cells.Item(107, 8).Value = 2648.7  # H107
$ws.Cells.Item(107, 9).Value = 2366.3684  # I107
$ws.Cells.Item(107, 11).Value = 2366.3684  # K107
$ws.Cells.Item(107, 13).Value = -446.3683999999998  # M107

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1926.125  # H16
$ws.Cells.Item(16, 9).Value = 2130  # I16
$ws.Cells.Item(16, 11).Value = 2130  # K16
$ws.Cells.Item(16, 13).Value = -1843  # M16

$ws.Cells.Item(99, 8).Value = 3506.9092  # H99
$ws.Cells.Item(99, 9).Value = 2013  # I99
$ws.Cells.Item(99, 11).Value = 2013  # K99
$ws.Cells.Item(99, 13).Value = -515  # M99

$ws.Cells.Item(107, 8).Value = 521.5333000000001  # H107
$ws.Cells.Item(107, 9).Value = 445.25  # I107
$ws.Cells.Item(107, 11).Value = 445.25  # K107
$ws.Cells.Item(107, 13).Value = 1474.75  # M107

$ws.Cells.Item(113, 8).Value = 1926.125  # H113
$ws.Cells.Item(113, 9).Value = 2130  # I113
$ws.Cells.Item(113, 11).Value = 2130  # K113
$ws.Cells.Item(113, 13).Value = 40  # M113

$ws.Cells.Item(126, 8).Value = 3506.9092  # H126
$ws.Cells.Item(126, 9).Value = 2013  # I126
$ws.Cells.Item(126, 11).Value = 6039  # K126
$ws.Cells.Item(126, 13).Value = -3569  # M126

$ws.Cells.Item(134, 8).Value = 1740.8438  # H134
$ws.Cells.Item(134, 9).Value = 1674.4839  # I134
$ws.Cells.Item(134, 11).Value = 5023.4517  # K134
$ws.Cells.Item(134, 13).Value = -2488.4517  # M134

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 1446.5333  # H14
$ws.Cells.Item(14, 9).Value = 1446.5333  # I14
$ws.Cells.Item(14, 11).Value = 4339.5999  # K14
$ws.Cells.Item(14, 13).Value = -4166.5999  # M14

$ws.Cells.Item(68, 13).ClearContents()  # M68
$ws.Cells.Item(68, 8).Value = 20840516  # H68
$ws.Cells.Item(68, 9).Value = 0  # I68
$ws.Cells.Item(68, 11).Value = 0  # K68

$ws.Cells.Item(71, 13).ClearContents()  # M71
$ws.Cells.Item(71, 8).Value = 20840516  # H71
$ws.Cells.Item(71, 9).Value = 0  # I71
$ws.Cells.Item(71, 11).Value = 0  # K71

$ws.Cells.Item(132, 8).Value = 4594.952  # H132
$ws.Cells.Item(132, 9).Value = 2621.3572  # I132
$ws.Cells.Item(132, 10).Value = 8542.143  # J132
$ws.Cells.Item(132, 11).Value = 23592.2148  # K132
$ws.Cells.Item(132, 12).Value = 76879.287  # L132
$ws.Cells.Item(132, 13).Value = -21062.2148  # M132
$ws.Cells.Item(132, 14).Value = -81939.287  # N132

$ws.Cells.Item(134, 8).Value = 16452.777  # H134
$ws.Cells.Item(134, 10).Value = 18384.375  # J134
$ws.Cells.Item(134, 12).Value = 55153.125  # L134
$ws.Cells.Item(134, 14).Value = -65293.125  # N134

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 132.07692  # H2
$ws.Cells.Item(2, 9).Value = 166.6  # I2
$ws.Cells.Item(2, 11).Value = 166.6  # K2
$ws.Cells.Item(2, 13).Value = -53.59999999999999  # M2

$ws.Cells.Item(70, 14).ClearContents()  # N70
$ws.Cells.Item(70, 8).Value = 4600.8  # H70
$ws.Cells.Item(70, 9).Value = 4600.8  # I70
$ws.Cells.Item(70, 10).Value = 0  # J70
$ws.Cells.Item(70, 11).Value = 4600.8  # K70
$ws.Cells.Item(70, 12).Value = 0  # L70
$ws.Cells.Item(70, 13).Value = -4330.8  # M70

$ws.Cells.Item(73, 14).ClearContents()  # N73
$ws.Cells.Item(73, 8).Value = 4600.8  # H73
$ws.Cells.Item(73, 9).Value = 4600.8  # I73
$ws.Cells.Item(73, 10).Value = 0  # J73
$ws.Cells.Item(73, 11).Value = 4600.8  # K73
$ws.Cells.Item(73, 12).Value = 0  # L73
$ws.Cells.Item(73, 13).Value = -3664.8  # M73

$ws.Cells.Item(102, 8).Value = 1438.8  # H102
$ws.Cells.Item(102, 9).Value = 1438.8  # I102
$ws.Cells.Item(102, 11).Value = 1438.8  # K102
$ws.Cells.Item(102, 13).Value = 183.2  # M102

$ws.Cells.Item(107, 8).Value = 932.25  # H107
$ws.Cells.Item(107, 9).Value = 327.8  # I107
$ws.Cells.Item(107, 10).Value = 9999  # J107
$ws.Cells.Item(107, 11).Value = 327.8  # K107
$ws.Cells.Item(107, 12).Value = 9999  # L107
$ws.Cells.Item(107, 13).Value = 1592.2  # M107
$ws.Cells.Item(107, 14).Value = -13839  # N107

$ws.Cells.Item(122, 8).Value = 3425.3044  # H122
$ws.Cells.Item(122, 9).Value = 2898.8235  # I122
$ws.Cells.Item(122, 11).Value = 8696.470499999999  # K122
$ws.Cells.Item(122, 13).Value = -6246.470499999999  # M122

$ws.Cells.Item(126, 8).Value = 4113.8335  # H126
$ws.Cells.Item(126, 10).Value = 4499.5  # J126
$ws.Cells.Item(126, 12).Value = 13498.5  # L126
$ws.Cells.Item(126, 14).Value = -18438.5  # N126

$ws.Cells.Item(132, 8).Value = 2983  # H132
$ws.Cells.Item(132, 9).Value = 2983  # I132
$ws.Cells.Item(132, 11).Value = 8949  # K132
$ws.Cells.Item(132, 13).Value = -6419  # M132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 3466666.8  # H2
$ws.Cells.Item(2, 9).Value = 200000  # I2
$ws.Cells.Item(2, 10).Value = 10000000  # J2
$ws.Cells.Item(2, 11).Value = 200000  # K2
$ws.Cells.Item(2, 12).Value = 10000000  # L2
$ws.Cells.Item(2, 13).Value = -199888  # M2
$ws.Cells.Item(2, 14).Value = -10000224  # N2

$ws.Cells.Item(19, 8).Value = 260  # H19
$ws.Cells.Item(19, 9).Value = 146.66667  # I19
$ws.Cells.Item(19, 10).Value = 600  # J19
$ws.Cells.Item(19, 11).Value = 146.66667  # K19
$ws.Cells.Item(19, 12).Value = 600  # L19
$ws.Cells.Item(19, 13).Value = 23.33332999999999  # M19
$ws.Cells.Item(19, 14).Value = -940  # N19

$ws.Cells.Item(22, 8).Value = 1063  # H22
$ws.Cells.Item(22, 9).Value = 864.6667  # I22
$ws.Cells.Item(22, 10).Value = 1241.5  # J22
$ws.Cells.Item(22, 11).Value = 864.6667  # K22
$ws.Cells.Item(22, 12).Value = 1241.5  # L22
$ws.Cells.Item(22, 13).Value = -569.6667  # M22
$ws.Cells.Item(22, 14).Value = -1831.5  # N22

$ws.Cells.Item(27, 8).Value = 1063  # H27
$ws.Cells.Item(27, 9).Value = 864.6667  # I27
$ws.Cells.Item(27, 10).Value = 1241.5  # J27
$ws.Cells.Item(27, 11).Value = 864.6667  # K27
$ws.Cells.Item(27, 12).Value = 1241.5  # L27
$ws.Cells.Item(27, 13).Value = -757.6667  # M27
$ws.Cells.Item(27, 14).Value = -1455.5  # N27

$ws.Cells.Item(40, 8).Value = 3610.0667  # H40
$ws.Cells.Item(40, 9).Value = 2544  # I40
$ws.Cells.Item(40, 11).Value = 2544  # K40
$ws.Cells.Item(40, 13).Value = -2408  # M40

$ws.Cells.Item(61, 8).Value = 3073.0605  # H61
$ws.Cells.Item(61, 9).Value = 2795.3462  # I61
$ws.Cells.Item(61, 11).Value = 2795.3462  # K61
$ws.Cells.Item(61, 13).Value = -2593.3462  # M61

$ws.Cells.Item(113, 8).Value = 3073.0605  # H113
$ws.Cells.Item(113, 9).Value = 2795.3462  # I113
$ws.Cells.Item(113, 11).Value = 2795.3462  # K113
$ws.Cells.Item(113, 13).Value = -625.3462  # M113

$ws.Cells.Item(139, 8).Value = 89600  # H139
$ws.Cells.Item(139, 9).Value = 89600  # I139
$ws.Cells.Item(139, 11).Value = 89600  # K139
$ws.Cells.Item(139, 13).Value = -84460  # M139

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(40, 14).ClearContents()  # N40
$ws.Cells.Item(40, 8).Value = 0  # H40
$ws.Cells.Item(40, 10).Value = 0  # J40
$ws.Cells.Item(40, 12).Value = 0  # L40

$ws.Cells.Item(69, 8).Value = 21722  # H69
$ws.Cells.Item(69, 10).Value = 21722  # J69
$ws.Cells.Item(69, 12).Value = 21722  # L69
$ws.Cells.Item(69, 14).Value = -23220  # N69

$ws.Cells.Item(72, 8).Value = 21722  # H72
$ws.Cells.Item(72, 10).Value = 21722  # J72
$ws.Cells.Item(72, 12).Value = 65166  # L72
$ws.Cells.Item(72, 14).Value = -72654  # N72

$ws.Cells.Item(96, 8).Value = 5623.75  # H96
$ws.Cells.Item(96, 9).Value = 6831.6665  # I96
$ws.Cells.Item(96, 11).Value = 6831.6665  # K96
$ws.Cells.Item(96, 13).Value = -5458.6665  # M96
